$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (insert shared strings in the target canonical order:
# 30=Nb Iter, 31=Nb Ctr ajoutee, 32=Solution, 33=Tps + CPU)
$ws.Range("D1").Value = "Nb Iter"
$ws.Range("E1").Value = "Nb Ctr ajoutée"
$ws.Range("B1").Value = "Solution"
$ws.Range("C1").Value = "Tps + CPU"

# Column C timing values, written in shared-string first-seen order
$ws.Range("C2").Value = "0.137014 seconds (5.89 k allocations: 393.406 KiB)"
$ws.Range("C3").Value = "0.165501 seconds (26.36 k allocations: 2.154 MiB)"
$ws.Range("C4").Value = "0.206548 seconds (57.40 k allocations: 5.059 MiB)"
$ws.Range("C5").Value = "1.062222 seconds (123.52 k allocations: 11.635 MiB)"
$ws.Range("C7").Value = "3.375432 seconds (311.61 k allocations: 31.073 MiB, 0.18% gc time)"
$ws.Range("C6").Value = "0.492398 seconds (195.61 k allocations: 19.725 MiB, 1.01% gc time)"
$ws.Range("C8").Value = "17.763255 seconds (667.01 k allocations: 72.797 MiB, 0.10% gc time)"
$ws.Range("C9").Value = "6.156659 seconds (685.97 k allocations: 73.051 MiB, 0.29% gc time)"
$ws.Range("C10").Value = "7.125538 seconds (911.33 k allocations: 96.495 MiB, 0.28% gc time)"
$ws.Range("C11").Value = "41.788991 seconds (1.41 M allocations: 157.975 MiB, 0.38% gc time)"
$ws.Range("C13").Value = "35.971777 seconds (2.14 M allocations: 236.806 MiB, 0.13% gc time)"
$ws.Range("C16").Value = "105.200705 seconds (3.88 M allocations: 445.682 MiB, 0.08% gc time)"
$ws.Range("C12").Value = "12.828082 seconds (1.42 M allocations: 155.523 MiB, 0.27% gc time)"
$ws.Range("C14").Value = "45.006563 seconds (2.75 M allocations: 323.459 MiB, 0.14% gc time)"
$ws.Range("C15").Value = "74.752989 seconds (3.30 M allocations: 383.730 MiB, 0.25% gc time)"
$ws.Range("C18").Value = "  0.130585 seconds (5.89 k allocations: 394.094 KiB)"
$ws.Range("C19").Value = "  0.145352 seconds (13.47 k allocations: 937.688 KiB)"
$ws.Range("C20").Value = "  0.131483 seconds (6.80 k allocations: 856.938 KiB)"
$ws.Range("C21").Value = "  0.155768 seconds (49.13 k allocations: 3.593 MiB)"
$ws.Range("C22").Value = "  0.480462 seconds (110.15 k allocations: 9.687 MiB, 1.22% gc time)"
$ws.Range("C23").Value = "  0.478024 seconds (154.46 k allocations: 13.213 MiB, 1.35% gc time)"
$ws.Range("C24").Value = "  0.384103 seconds (175.00 k allocations: 14.856 MiB)"
$ws.Range("C25").Value = "  1.156912 seconds (296.54 k allocations: 26.947 MiB, 0.99% gc time)"
$ws.Range("C26").Value = "  0.976702 seconds (302.92 k allocations: 25.066 MiB, 0.66% gc time)"
$ws.Range("C27").Value = "  3.671390 seconds (521.89 k allocations: 50.510 MiB, 0.50% gc time)"
$ws.Range("C28").Value = "  1.016204 seconds (399.29 k allocations: 31.878 MiB, 0.65% gc time)"
$ws.Range("C29").Value = "  2.732261 seconds (655.98 k allocations: 58.963 MiB, 0.89% gc time)"
$ws.Range("C30").Value = "  0.819629 seconds (487.95 k allocations: 38.457 MiB, 1.75% gc time)"
$ws.Range("C31").Value = "  9.981334 seconds (1.30 M allocations: 134.869 MiB, 1.38% gc time)"
$ws.Range("C32").Value = "  3.296910 seconds (1.02 M allocations: 94.883 MiB, 1.07% gc time)"
